$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing Key column (A) values for rows 2-4
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# New rows' numeric Key values
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6

# New text values, entered in the order that reproduces the authored
# shared-string table layout
$ws.Range("C5").Value = "user@abv.bg"
$ws.Range("E5").Value = "password"
$ws.Range("F5").Value = "password"

$ws.Range("D6").Value = "User"

$ws.Range("F7").Value = "password1"
$ws.Range("E7").Value = "password"
$ws.Range("D7").Value = "User"
$ws.Range("C7").Value = "user@abv.bg"

$ws.Range("B7").Value = "RegisterMIsmatchPassword"
$ws.Range("B5").Value = "RegisterWithoutFullName"
$ws.Range("B6").Value = "RegisterWithoutPassword"
$ws.Range("C6").Value = "user@abv.bg"

# Update the selection to match the authored state
$ws.Range("B6").Select()
